# Update "想去人数" (want-to-go count) values for two events whose data
# appears on both the "展览" sheet and the "全部类型" sheet.
#   id=89145 event: 1247 -> 1250
#   id=88276 event: 2736 -> 2741

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1250
$wsExhibit.Range("F4").Value = 2741

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1250
$wsAll.Range("F6").Value = 2741
